$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "pages" (sheet1)
# ---------------------------------------------------------------------------
$pages = $wb.Worksheets.Item("pages")

# Insert a brand-new row 6 (shifts old rows 6-8 down to 7-9, keeping the
# formatting of the row above, same as Excel's normal "Insert" behaviour).
$pages.Rows("6:6").Insert()
$pages.Rows("6:6").RowHeight = 24

# New row 6: an "offline" / page-refresh action
$pages.Range("C6").Value = "page-refresh"
# Row 5: the "other" layer entry becomes "offline"
$pages.Range("B5").Value = "offline"
$pages.Range("F6").Value = "somethingDifferent"
# Row 4: fix the "somthingCustom" typo -> "somethingCustom"
$pages.Range("F4").Value = "somethingCustom"

# finish filling in row 5 (offline / page-load / actions / somethingGlobal)
$pages.Range("C5").Value = "page-load"
$pages.Range("D5").Value = "actions"
$pages.Range("F5").Value = "somethingGlobal"

# finish filling in the new row 6 (blank / page-refresh / actions / somethingDifferent)
$pages.Range("B6").Value = ""
$pages.Range("D6").Value = "actions"

# Row 8 (previously row 7, shifted down by the insert): "other"/"page-other" -> "offline"/"page-refresh"
$pages.Range("B8").Value = "offline"
$pages.Range("C8").Value = "page-refresh"
$pages.Range("D8").Value = "actions"
$pages.Range("F8").Value = "somethingSpecific"

# ---------------------------------------------------------------------------
# Sheet "events" (sheet2)
# ---------------------------------------------------------------------------
$events = $wb.Worksheets.Item("events")

# Fix the same "somthingCustom" typo -> "somethingCustom"
$events.Range("F4").Value = "somethingCustom"

# Narrow columns A:F from 29 to 23 characters wide
$events.Range("A1:F1").ColumnWidth = 22.1666666666667

# ---------------------------------------------------------------------------
# View state: make "pages" the selected / active sheet with its own
# selection, then restore the selection on "events" too. Selecting a range
# on a sheet activates that sheet, so "events" is touched first (just for
# its own selection) and "pages" is touched last so it ends up the active
# tab, matching the target workbook view.
# ---------------------------------------------------------------------------
$events.Range("F18").Select()
$pages.Range("C19").Select()
